# Update the answer table's division-problem results.
# Each row below is: (WordRow, WordCol, NewText) addressing the 1st table
# via Table.Cell(row, col) (1-based). Only the text content changes; run
# formatting (font/size) is preserved by assigning to Cell.Range.Text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @(1, 1, "90÷4=22, 2"),
    @(1, 2, "45÷4=11, 1"),
    @(1, 3, "19÷3=6, 1"),
    @(1, 4, "18÷4=4, 2"),
    @(1, 5, "76÷8=9, 4"),

    @(5, 1, "82÷8=10, 2"),
    @(5, 2, "31÷9=3, 4"),
    @(5, 3, "96÷7=13, 5"),
    @(5, 4, "16÷5=3, 1"),
    @(5, 5, "29÷2=14, 1"),

    @(9, 1, "21÷2=10, 1"),
    @(9, 2, "70÷4=17, 2"),
    @(9, 3, "81÷5=16, 1"),
    @(9, 4, "47÷4=11, 3"),
    @(9, 5, "98÷6=16, 2"),

    @(13, 1, "90÷6=15, 0"),
    @(13, 2, "19÷2=9, 1"),
    @(13, 3, "32÷9=3, 5"),
    @(13, 4, "63÷7=9, 0"),
    @(13, 5, "31÷4=7, 3"),

    @(17, 1, "83÷6=13, 5"),
    @(17, 2, "33÷2=16, 1"),
    @(17, 3, "69÷4=17, 1"),
    @(17, 4, "39÷4=9, 3"),
    @(17, 5, "20÷9=2, 2")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $text = $u[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $text
}
